$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated statistics after correcting selection scopes (each row's N increased by 1)
$ws.Range("B2").Value = 0.1504689633014369
$ws.Range("C2").Value = 0.2334218385963792
$ws.Range("D2").Value = 0.2642580679566157
$ws.Range("E2").Value = 0.5140603738439832
$ws.Range("F2").Value = 0.5031130469799929
$ws.Range("G2").Value = 22

$ws.Range("B3").Value = 0.406982015469362
$ws.Range("C3").Value = 2.696787301862778
$ws.Range("D3").Value = 27.56079346563211
$ws.Range("E3").Value = 5.249837470401546
$ws.Range("F3").Value = 5.36329349000709
$ws.Range("G3").Value = 21

$ws.Range("B4").Value = -0.9836666976115556
$ws.Range("C4").Value = 2.507506085034803
$ws.Range("D4").Value = 19.95067858449099
$ws.Range("E4").Value = 4.466618249245283
$ws.Range("F4").Value = 4.470144047190493
$ws.Range("G4").Value = 20

$ws.Range("B5").Value = -0.1790030931086928
$ws.Range("C5").Value = 1.392189940396086
$ws.Range("D5").Value = 6.735992462938008
$ws.Range("E5").Value = 2.595379059586096
$ws.Range("F5").Value = 2.660148875912306
$ws.Range("G5").Value = 19

$ws.Range("B6").Value = -0.08048131822523331
$ws.Range("C6").Value = 1.576718745090598
$ws.Range("D6").Value = 10.20882528216252
$ws.Range("E6").Value = 3.195125237320522
$ws.Range("F6").Value = 3.28671358039524
$ws.Range("G6").Value = 18

$ws.Range("B7").Value = -0.3841978569163946
$ws.Range("C7").Value = 2.012713761047529
$ws.Range("D7").Value = 14.27256267236992
$ws.Range("E7").Value = 3.777904534576003
$ws.Range("F7").Value = 3.873985589358225
$ws.Range("G7").Value = 17

$ws.Range("B8").Value = -0.2198684980862554
$ws.Range("C8").Value = 1.996418041492534
$ws.Range("D8").Value = 11.77248202587504
$ws.Range("E8").Value = 3.431105073569599
$ws.Range("F8").Value = 3.544231712448563
$ws.Range("G8").Value = 15
